$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Octubre de 2020 a las 11:28"

# Row 7
$ws.Cells.Item(7, 2).Value = 1260112
$ws.Cells.Item(7, 3).Value = 11493
$ws.Cells.Item(7, 4).Value = 1002329
$ws.Cells.Item(7, 5).Value = 235727
$ws.Cells.Item(7, 7).Value = 191
$ws.Cells.Item(7, 8).Value = 22056

# Row 24
$ws.Cells.Item(24, 1).Value = "Indonesia"
$ws.Cells.Item(24, 2).Value = 320564
$ws.Cells.Item(24, 3).Value = 4850
$ws.Cells.Item(24, 4).Value = 244060
$ws.Cells.Item(24, 5).Value = 64924
$ws.Cells.Item(24, 7).Value = 108
$ws.Cells.Item(24, 8).Value = 11580

# Row 25
$ws.Cells.Item(25, 1).Value = "Pakistan"
$ws.Cells.Item(25, 2).Value = 316934
$ws.Cells.Item(25, 3).Value = 583
$ws.Cells.Item(25, 4).Value = 302375
$ws.Cells.Item(25, 5).Value = 8015
$ws.Cells.Item(25, 7).Value = 9
$ws.Cells.Item(25, 8).Value = 6544

# Row 27
$ws.Cells.Item(27, 2).Value = 283532
$ws.Cells.Item(27, 3).Value = 2051
$ws.Cells.Item(27, 4).Value = 220046
$ws.Cells.Item(27, 5).Value = 61640
$ws.Cells.Item(27, 7).Value = 22
$ws.Cells.Item(27, 8).Value = 1846

# Row 39
$ws.Cells.Item(39, 1).Value = "Polonia"
$ws.Cells.Item(39, 2).Value = 111599
$ws.Cells.Item(39, 3).Value = 4280
$ws.Cells.Item(39, 4).Value = 76490
$ws.Cells.Item(39, 5).Value = 32242
$ws.Cells.Item(39, 7).Value = 76
$ws.Cells.Item(39, 8).Value = 2867

# Row 40
$ws.Cells.Item(40, 1).Value = "Kuwait"
$ws.Cells.Item(40, 2).Value = 108743
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 100776
$ws.Cells.Item(40, 5).Value = 7328
$ws.Cells.Item(40, 8).Value = 639

# Row 41
$ws.Cells.Item(41, 1).Value = "Kazajistan"
$ws.Cells.Item(41, 2).Value = 108454
$ws.Cells.Item(41, 3).Value = 92
$ws.Cells.Item(41, 4).Value = 103604
$ws.Cells.Item(41, 5).Value = 3104
$ws.Cells.Item(41, 8).Value = 1746

# Row 65
$ws.Cells.Item(65, 2).Value = 52057
$ws.Cells.Item(65, 3).Value = 1209
$ws.Cells.Item(65, 4).Value = 41289
$ws.Cells.Item(65, 5).Value = 9930
$ws.Cells.Item(65, 7).Value = 8
$ws.Cells.Item(65, 8).Value = 838

# Row 91
$ws.Cells.Item(91, 2).Value = 18989
$ws.Cells.Item(91, 3).Value = 542
$ws.Cells.Item(91, 4).Value = 16473
$ws.Cells.Item(91, 5).Value = 2206
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 310

# Row 93
$ws.Cells.Item(93, 1).Value = "Eslovaquia"
$ws.Cells.Item(93, 2).Value = 15726
$ws.Cells.Item(93, 3).Value = 1037
$ws.Cells.Item(93, 4).Value = 5353
$ws.Cells.Item(93, 5).Value = 10316
$ws.Cells.Item(93, 7).Value = 2
$ws.Cells.Item(93, 8).Value = 57

# Row 94
$ws.Cells.Item(94, 1).Value = "Zambia"
$ws.Cells.Item(94, 2).Value = 15224
$ws.Cells.Item(94, 4).Value = 14342
$ws.Cells.Item(94, 5).Value = 547
$ws.Cells.Item(94, 8).Value = 335

# Row 95
$ws.Cells.Item(95, 1).Value = "Senegal"
$ws.Cells.Item(95, 2).Value = 15174
$ws.Cells.Item(95, 4).Value = 12998
$ws.Cells.Item(95, 5).Value = 1863
$ws.Cells.Item(95, 8).Value = 313

# Row 96
$ws.Cells.Item(96, 1).Value = "Noruega"
$ws.Cells.Item(96, 2).Value = 15013
$ws.Cells.Item(96, 4).Value = 11863
$ws.Cells.Item(96, 5).Value = 2875
$ws.Cells.Item(96, 8).Value = 275

# Row 97
$ws.Cells.Item(97, 1).Value = "Albania"
$ws.Cells.Item(97, 2).Value = 14730
$ws.Cells.Item(97, 4).Value = 9115
$ws.Cells.Item(97, 5).Value = 5208
$ws.Cells.Item(97, 8).Value = 407

# Row 116
$ws.Cells.Item(116, 1).Value = "Eslovenia"
$ws.Cells.Item(116, 2).Value = 7507
$ws.Cells.Item(116, 3).Value = 387
$ws.Cells.Item(116, 4).Value = 4655
$ws.Cells.Item(116, 5).Value = 2692
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 160

# Row 117
$ws.Cells.Item(117, 1).Value = "Jamaica"
$ws.Cells.Item(117, 2).Value = 7191
$ws.Cells.Item(117, 4).Value = 2700
$ws.Cells.Item(117, 5).Value = 4365
$ws.Cells.Item(117, 8).Value = 126

# Row 153
$ws.Cells.Item(153, 1).Value = "Letonia"
$ws.Cells.Item(153, 2).Value = 2370
$ws.Cells.Item(153, 3).Value = 109
$ws.Cells.Item(153, 4).Value = 1322
$ws.Cells.Item(153, 5).Value = 1008
$ws.Cells.Item(153, 8).Value = 40

# Row 154
$ws.Cells.Item(154, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(154, 2).Value = 2358
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 1857
$ws.Cells.Item(154, 5).Value = 491
$ws.Cells.Item(154, 8).Value = 10

# Row 155
$ws.Cells.Item(155, 1).Value = "Belice"
$ws.Cells.Item(155, 2).Value = 2310
$ws.Cells.Item(155, 3).Value = 67
$ws.Cells.Item(155, 4).Value = 1427
$ws.Cells.Item(155, 5).Value = 849
$ws.Cells.Item(155, 8).Value = 34

# Row 156
$ws.Cells.Item(156, 1).Value = "Sierra Leona"
$ws.Cells.Item(156, 2).Value = 2287
$ws.Cells.Item(156, 4).Value = 1716
$ws.Cells.Item(156, 5).Value = 499
$ws.Cells.Item(156, 8).Value = 72

# Row 175
$ws.Cells.Item(175, 2).Value = 524
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(175, 5).Value = 32

# Row 196
$ws.Cells.Item(196, 2).Value = 135
$ws.Cells.Item(196, 3).Value = 4
$ws.Cells.Item(196, 4).Value = 117
$ws.Cells.Item(196, 5).Value = 17

# Row 207
$ws.Cells.Item(207, 1).Value = "Santa Lucia"

# Row 208
$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"

# Row 215
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0

# Row 216
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 8).Value = 1
